# Update "想去人数" (interest count) figures for a few events that changed
# between the previous and the newly generated data pull.
#
# Sheet "展览" (sheet1) and sheet "全部类型" (sheet4) both list the same
# events (全部类型 aggregates all the other sheets), so each value needs to
# be updated in both places.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 蜀山·银泰百货高新店-2024漫趣地带嘉年华（免费） : 269 -> 268
$wsExhibit.Range("F2").Value = 268
$wsAll.Range("F2").Value = 268

# 合肥·第十五届次元之门动漫游戏博览会 : 6679 -> 6680
$wsExhibit.Range("F5").Value = 6680
$wsAll.Range("F5").Value = 6680

# 合肥·第九届环形宇宙动漫游戏嘉年华 : 121 -> 123
# (row 12 on "展览", row 14 on "全部类型")
$wsExhibit.Range("F12").Value = 123
$wsAll.Range("F14").Value = 123
